# Auto-generated Excel COM-interop script applying the Kujata_Profits profit-recalculation update.
# For each affected worksheet, a set of (cell, value) pairs is applied. A value of $null
# means the cell had its content cleared (column removed from the row).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (60 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    ,("H33", 627)
    ,("I33", 643.3333)
    ,("K33", 643.3333)
    ,("M33", -414.3333)
    ,("H80", 894.93335)
    ,("I80", 1573)
    ,("J80", 648.36365)
    ,("K80", 4719)
    ,("L80", 1945.09095)
    ,("M80", -3721)
    ,("N80", -3941.09095)
    ,("H83", 894.93335)
    ,("I83", 1573)
    ,("J83", 648.36365)
    ,("K83", 14157)
    ,("L83", 5835.27285)
    ,("M83", -9165)
    ,("N83", -15819.27285)
    ,("H88", 1123280.6)
    ,("I88", 380)
    ,("J88", 1764938.1)
    ,("K88", 380)
    ,("L88", 1764938.1)
    ,("M88", 26)
    ,("N88", -1765750.1)
    ,("H91", 1123280.6)
    ,("I91", 380)
    ,("J91", 1764938.1)
    ,("K91", 380)
    ,("L91", 1764938.1)
    ,("M91", 1024)
    ,("N91", -1767746.1)
    ,("I135", 381.6316)
    ,("J135", 125001700)
    ,("K135", 3434.6844)
    ,("L135", 1125015300)
    ,("M135", -899.6844000000001)
    ,("N135", -1125020370)
    ,("H137", 1347.25)
    ,("I137", 933.3333)
    ,("J137", 1657.6875)
    ,("K137", 2799.9999)
    ,("L137", 4973.0625)
    ,("M137", -249.9998999999998)
    ,("N137", -10073.0625)
    ,("H138", 1516.92)
    ,("I138", 656.43243)
    ,("J138", 2022.2858)
    ,("K138", 1969.29729)
    ,("L138", 6066.857400000001)
    ,("M138", 3170.70271)
    ,("N138", -16346.8574)
    ,("H140", 39824)
    ,("J140", 39824)
    ,("H141", 795)
    ,("I141", 795)
    ,("K141", 2385)
    ,("M141", 2795)
    ,("L140", 39824)
    ,("N140", -50184)
)
foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    if ($null -eq $newVal) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = $newVal
    }
}

# ---- Sheet: ARM (75 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    ,("H2", 8280.77)
    ,("I2", 584.5714)
    ,("J2", 17259.666)
    ,("K2", 584.5714)
    ,("L2", 17259.666)
    ,("M2", -471.5714)
    ,("N2", -17485.666)
    ,("H32", 3017.1448)
    ,("I32", 2847.2698)
    ,("J32", 3840.3845)
    ,("K32", 2847.2698)
    ,("L32", 3840.3845)
    ,("M32", -2560.2698)
    ,("N32", -4414.3845)
    ,("H45", 1113.2778)
    ,("I45", 1074.2142)
    ,("K45", 1074.2142)
    ,("M45", -697.2141999999999)
    ,("H61", 1089.8823)
    ,("I61", 1015.4815)
    ,("J61", 1376.8572)
    ,("K61", 1015.4815)
    ,("L61", 1376.8572)
    ,("M61", -803.4815)
    ,("N61", -1800.8572)
    ,("H74", 1279.12)
    ,("I74", 796.4706)
    ,("J74", 2304.75)
    ,("K74", 796.4706)
    ,("L74", 2304.75)
    ,("M74", 77.52940000000001)
    ,("N74", -4052.75)
    ,("H77", 1279.12)
    ,("I77", 796.4706)
    ,("J77", 2304.75)
    ,("K77", 3982.353)
    ,("L77", 11523.75)
    ,("M77", 385.6469999999999)
    ,("N77", -20259.75)
    ,("H88", 2862.8125)
    ,("I88", 2201.25)
    ,("K88", 2201.25)
    ,("M88", -1795.25)
    ,("H91", 2862.8125)
    ,("I91", 2201.25)
    ,("K91", 2201.25)
    ,("M91", -797.25)
    ,("H116", 8280.77)
    ,("I116", 584.5714)
    ,("J116", 17259.666)
    ,("K116", 584.5714)
    ,("L116", 17259.666)
    ,("M116", 1709.4286)
    ,("N116", -21847.666)
    ,("H122", 810.6070999999999)
    ,("I122", 824.0417)
    ,("J122", 730)
    ,("K122", 2472.1251)
    ,("L122", 2190)
    ,("M122", -22.1251000000002)
    ,("N122", -7090)
    ,("H132", 1532.4688)
    ,("I132", 1186.7037)
    ,("J132", 3399.6)
    ,("K132", 3560.1111)
    ,("L132", 10198.8)
    ,("M132", -1030.1111)
    ,("N132", -15258.8)
    ,("H136", 1089.8823)
    ,("I136", 1015.4815)
    ,("J136", 1376.8572)
    ,("K136", 3046.4445)
    ,("L136", 4130.571599999999)
    ,("M136", -496.4445000000001)
    ,("N136", -9230.571599999999)
)
foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    if ($null -eq $newVal) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = $newVal
    }
}

# ---- Sheet: BSM (40 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    ,("H3", 8280.77)
    ,("I3", 584.5714)
    ,("J3", 17259.666)
    ,("K3", 584.5714)
    ,("L3", 17259.666)
    ,("M3", -470.5714)
    ,("N3", -17487.666)
    ,("H86", 4079.72)
    ,("I86", 4045.182)
    ,("K86", 4045.182)
    ,("M86", -2922.182)
    ,("H89", 4079.72)
    ,("I89", 4045.182)
    ,("K89", 20225.91)
    ,("M89", -14609.91)
    ,("H94", 25001014)
    ,("J94", 1210)
    ,("L94", 1210)
    ,("N94", -2112)
    ,("H99", 33334480)
    ,("I99", 50001100)
    ,("J99", 1242.2)
    ,("K99", 50001100)
    ,("L99", 1242.2)
    ,("M99", -49999602)
    ,("N99", -4238.2)
    ,("H105", 333336060)
    ,("I105", 333336060)
    ,("J105", 0)
    ,("K105", 333336060)
    ,("L105", 0)
    ,("H134", 1528.7954)
    ,("I134", 1016.45715)
    ,("J134", 3521.2222)
    ,("K134", 3049.37145)
    ,("L134", 10563.6666)
    ,("M134", -514.3714499999996)
    ,("N134", -15633.6666)
    ,("M105", -333334313)
    ,("N105", $null)
)
foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    if ($null -eq $newVal) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = $newVal
    }
}

# ---- Sheet: CRP (42 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    ,("H58", 1337.125)
    ,("I58", 1234.5454)
    ,("J58", 1562.8)
    ,("K58", 1234.5454)
    ,("L58", 1562.8)
    ,("M58", -1031.5454)
    ,("N58", -1968.8)
    ,("H62", 7694634.5)
    ,("I62", 2419.96)
    ,("J62", 200000000)
    ,("K62", 2419.96)
    ,("L62", 200000000)
    ,("M62", -1795.96)
    ,("N62", -200001248)
    ,("H65", 7694634.5)
    ,("I65", 2419.96)
    ,("J65", 200000000)
    ,("K65", 12099.8)
    ,("L65", 1000000000)
    ,("M65", -8979.799999999999)
    ,("N65", -1000006240)
    ,("H132", 2066.2222)
    ,("I132", 1099.5)
    ,("J132", 3999.6667)
    ,("K132", 3298.5)
    ,("L132", 11999.0001)
    ,("M132", -768.5)
    ,("N132", -17059.0001)
    ,("H134", 1477.5454)
    ,("I134", 1162.1666)
    ,("J134", 1856)
    ,("K134", 3486.4998)
    ,("L134", 5568)
    ,("M134", -951.4998000000001)
    ,("N134", -10638)
    ,("H136", 1337.125)
    ,("I136", 1234.5454)
    ,("J136", 1562.8)
    ,("K136", 3703.6362)
    ,("L136", 4688.4)
    ,("M136", -1153.6362)
    ,("N136", -9788.4)
)
foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    if ($null -eq $newVal) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = $newVal
    }
}

# ---- Sheet: CUL (4 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    ,("H131", 11238038)
    ,("J131", 2227.0122)
    ,("L131", 6681.0366)
    ,("N131", -16761.0366)
)
foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    if ($null -eq $newVal) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = $newVal
    }
}

# ---- Sheet: GSM (14 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    ,("H80", 2927)
    ,("I80", 2325.625)
    ,("J80", 3614.2856)
    ,("K80", 2325.625)
    ,("L80", 3614.2856)
    ,("M80", -1327.625)
    ,("N80", -5610.2856)
    ,("H83", 2927)
    ,("I83", 2325.625)
    ,("J83", 3614.2856)
    ,("K83", 11628.125)
    ,("L83", 18071.428)
    ,("M83", -6636.125)
    ,("N83", -28055.428)
)
foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    if ($null -eq $newVal) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = $newVal
    }
}

# ---- Sheet: LTW (30 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    ,("H16", 805.5)
    ,("I16", 805.5)
    ,("K16", 805.5)
    ,("M16", -635.5)
    ,("H33", 5008.5)
    ,("J33", 5008.5)
    ,("L33", 5008.5)
    ,("N33", -5588.5)
    ,("H100", 2500)
    ,("J100", 2500)
    ,("H105", 17000)
    ,("J105", 17000)
    ,("H132", 24364.727)
    ,("I132", 1029.8)
    ,("J132", 55068.58)
    ,("K132", 3089.4)
    ,("L132", 165205.74)
    ,("M132", -559.3999999999996)
    ,("N132", -170265.74)
    ,("H136", 2032.6666)
    ,("I136", 2321)
    ,("J136", 1802)
    ,("K136", 6963)
    ,("L136", 5406)
    ,("M136", -4413)
    ,("N136", -10506)
    ,("L100", 2500)
    ,("N100", -3582)
    ,("L105", 17000)
    ,("N105", -23988)
)
foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    if ($null -eq $newVal) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = $newVal
    }
}

# ---- Sheet: WVR (7 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    ,("H136", 1231.625)
    ,("I136", 704)
    ,("J136", 1307)
    ,("K136", 2112)
    ,("L136", 3921)
    ,("M136", 438)
    ,("N136", -9021)
)
foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    if ($null -eq $newVal) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = $newVal
    }
}
